$p = $ppt.ActivePresentation

# The deck's live theme (ppt/theme/theme2.xml, referenced by the slide
# master / presentation) is switched from the custom "Integral" / "Red
# Violet" colour scheme to the stock "Office Theme" colour scheme. The
# font scheme and format scheme are left untouched (they were already
# identical between the two themes).
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

$colors.Colors(1).RGB  = 0        # dk1      000000
$colors.Colors(2).RGB  = 16777215 # lt1      FFFFFF
$colors.Colors(3).RGB  = 6968388  # dk2      44546A
$colors.Colors(4).RGB  = 15132391 # lt2      E7E6E6
$colors.Colors(5).RGB  = 13998939 # accent1  5B9BD5
$colors.Colors(6).RGB  = 3243501  # accent2  ED7D31
$colors.Colors(7).RGB  = 10855845 # accent3  A5A5A5
$colors.Colors(8).RGB  = 49407    # accent4  FFC000
$colors.Colors(9).RGB  = 12874308 # accent5  4472C4
$colors.Colors(10).RGB = 4697456  # accent6  70AD47
$colors.Colors(11).RGB = 12673797 # hlink    0563C1
$colors.Colors(12).RGB = 7491477  # folHlink 954F72
